$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = "9c5c7be0-5623-4c31-96cb-89cba41271c5"
$ws.Range("B8").Value = "Login with valid username and password"
$ws.Range("C8").Value = "PASSED"
$ws.Range("D8").Value = 45048.84917055556
$ws.Range("E8").Value = 45048.84922762732
$ws.Range("F8").Value = "PT4.9319502S"

# Row 9
$ws.Range("A9").Value = "9478cb78-099a-404d-8883-61f48c1d5fdd"
$ws.Range("B9").Value = "Create Country"
$ws.Range("C9").Value = "PASSED"
$ws.Range("D9").Value = 45048.849267743055
$ws.Range("E9").Value = 45048.8493509375
$ws.Range("F9").Value = "PT7.1877896S"
